# Update like-count / comment-count figures for the "环形宇宙动漫游戏嘉年华"
# and "《琳声雅集》" events across the two worksheets that list them
# ("展览" and "全部类型").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F6 1959 -> 1960, F11 65 -> 66
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 1960
$wsExhibit.Range("F11").Value = 66

# Sheet "全部类型" (sheet4): F7 1959 -> 1960, F12 65 -> 66
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 1960
$wsAll.Range("F12").Value = 66
